$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.641.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.185.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.180.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.700.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.61%  "

$ws.Range("E16").Value = "  -1.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.186.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.640.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.94%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.62%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0738"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("E42").Value = "  -4.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "398.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.805.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.112"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "
